# Auto-generated edit script applying NATMI LR-pair recompute (Fn1-Itga2b)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 7).Value = 40.7349555
$ws.Cells.Item(2, 8).Value = 81.469911
$ws.Cells.Item(2, 9).Value = 0.05567871843833241
$ws.Cells.Item(2, 10).Value = 0.03826666865920979
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 13).Value = 1.2502825
$ws.Cells.Item(2, 14).Value = 2.500565
$ws.Cells.Item(2, 15).Value = 0.1638687485091313
$ws.Cells.Item(2, 16).Value = 0.1262743888059313
$ws.Cells.Item(2, 17).Value = 50.93020199992875
$ws.Cells.Item(2, 18).Value = 203.720807999715
$ws.Cells.Item(2, 19).Value = 0.009124001909081829
$ws.Cells.Item(2, 20).Value = 0.004832100196580804

# Row 3
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 7).Value = 40.7349555
$ws.Cells.Item(3, 8).Value = 81.469911
$ws.Cells.Item(3, 9).Value = 0.05567871843833241
$ws.Cells.Item(3, 10).Value = 0.03826666865920979
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 2.596814
$ws.Cells.Item(3, 14).Value = 7.790442000000001
$ws.Cells.Item(3, 15).Value = 0.3403524085884521
$ws.Cells.Item(3, 16).Value = 0.3934044114342388
$ws.Cells.Item(3, 17).Value = 105.781102731777
$ws.Cells.Item(3, 18).Value = 634.686616390662
$ws.Cells.Item(3, 19).Value = 0.01895038592760469
$ws.Cells.Item(3, 20).Value = 0.01505427626142546

# Row 4
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 7).Value = 40.7349555
$ws.Cells.Item(4, 8).Value = 81.469911
$ws.Cells.Item(4, 9).Value = 0.05567871843833241
$ws.Cells.Item(4, 10).Value = 0.03826666865920979
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 0.2737596666666667
$ws.Cells.Item(4, 14).Value = 0.821279
$ws.Cells.Item(4, 15).Value = 0.03588041420154535
$ws.Cells.Item(4, 16).Value = 0.04147322855600493
$ws.Cells.Item(4, 17).Value = 11.1515878393615
$ws.Cells.Item(4, 18).Value = 66.909527036169
$ws.Cells.Item(4, 19).Value = 0.001997775479778587
$ws.Cells.Item(4, 20).Value = 0.001587042295380318

# Row 5
$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 7).Value = 40.7349555
$ws.Cells.Item(5, 8).Value = 81.469911
$ws.Cells.Item(5, 9).Value = 0.05567871843833241
$ws.Cells.Item(5, 10).Value = 0.03826666865920979
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 0.7929023333333333
$ws.Cells.Item(5, 14).Value = 2.378707
$ws.Cells.Item(5, 15).Value = 0.103922044060685
$ws.Cells.Item(5, 16).Value = 0.1201207617372036
$ws.Cells.Item(5, 17).Value = 32.29884126417949
$ws.Cells.Item(5, 18).Value = 193.793047585077
$ws.Cells.Item(5, 19).Value = 0.005786246230790855
$ws.Cells.Item(5, 20).Value = 0.004596621388489455

# Row 6
$ws.Cells.Item(6, 5).Value = 2
$ws.Cells.Item(6, 7).Value = 40.7349555
$ws.Cells.Item(6, 8).Value = 81.469911
$ws.Cells.Item(6, 9).Value = 0.05567871843833241
$ws.Cells.Item(6, 10).Value = 0.03826666865920979
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 0.879594
$ws.Cells.Item(6, 14).Value = 2.638782
$ws.Cells.Item(6, 15).Value = 0.1152843201245645
$ws.Cells.Item(6, 16).Value = 0.1332541182661091
$ws.Cells.Item(6, 17).Value = 35.830222448067
$ws.Cells.Item(6, 18).Value = 214.981334688402
$ws.Cells.Item(6, 19).Value = 0.006418883200570208
$ws.Cells.Item(6, 20).Value = 0.005099191191164351

# Row 7
$ws.Cells.Item(7, 5).Value = 2
$ws.Cells.Item(7, 7).Value = 39.5081721297966
$ws.Cells.Item(7, 8).Value = 81.469911
$ws.Cells.Item(7, 9).Value = 0.05567871843833241
$ws.Cells.Item(7, 10).Value = 0.03826666865920979
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 13).Value = 1.8364275
$ws.Cells.Item(7, 14).Value = 3.672855
$ws.Cells.Item(7, 15).Value = 0.2406920645156217
$ws.Cells.Item(7, 16).Value = 0.1854730912005122
$ws.Cells.Item(7, 17).Value = 74.80679249147624
$ws.Cells.Item(7, 18).Value = 299.227169965905
$ws.Cells.Item(7, 19).Value = 0.01340142569050624
$ws.Cells.Item(7, 20).Value = 0.007097437326169399

# Row 8
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 350.3919066666667
$ws.Cells.Item(8, 8).Value = 1051.17572
$ws.Cells.Item(8, 9).Value = 0.4789344206933965
$ws.Cells.Item(8, 10).Value = 0.4937404802104949
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 13).Value = 1.2502825
$ws.Cells.Item(8, 14).Value = 2.500565
$ws.Cells.Item(8, 15).Value = 0.1638687485091313
$ws.Cells.Item(8, 16).Value = 0.1262743888059313
$ws.Cells.Item(8, 17).Value = 438.0888690469666
$ws.Cells.Item(8, 18).Value = 2628.5332142818
$ws.Cells.Item(8, 19).Value = 0.07848238413697269
$ws.Cells.Item(8, 20).Value = 0.06234677736732729

# Row 9
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 350.3919066666667
$ws.Cells.Item(9, 8).Value = 1051.17572
$ws.Cells.Item(9, 9).Value = 0.4789344206933965
$ws.Cells.Item(9, 10).Value = 0.4937404802104949
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 2.596814
$ws.Cells.Item(9, 14).Value = 7.790442000000001
$ws.Cells.Item(9, 15).Value = 0.3403524085884521
$ws.Cells.Item(9, 16).Value = 0.3934044114342388
$ws.Cells.Item(9, 17).Value = 909.9026087186934
$ws.Cells.Item(9, 18).Value = 8189.12347846824
$ws.Cells.Item(9, 19).Value = 0.1630064836389125
$ws.Cells.Item(9, 20).Value = 0.1942396830184682

# Row 10
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 350.3919066666667
$ws.Cells.Item(10, 8).Value = 1051.17572
$ws.Cells.Item(10, 9).Value = 0.4789344206933965
$ws.Cells.Item(10, 10).Value = 0.4937404802104949
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 0.2737596666666667
$ws.Cells.Item(10, 14).Value = 0.821279
$ws.Cells.Item(10, 15).Value = 0.03588041420154535
$ws.Cells.Item(10, 16).Value = 0.04147322855600493
$ws.Cells.Item(10, 17).Value = 95.92317157176446
$ws.Cells.Item(10, 18).Value = 863.30854414588
$ws.Cells.Item(10, 19).Value = 0.01718436538985624
$ws.Cells.Item(10, 20).Value = 0.02047701178312148

# Row 11
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 350.3919066666667
$ws.Cells.Item(11, 8).Value = 1051.17572
$ws.Cells.Item(11, 9).Value = 0.4789344206933965
$ws.Cells.Item(11, 10).Value = 0.4937404802104949
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 0.7929023333333333
$ws.Cells.Item(11, 14).Value = 2.378707
$ws.Cells.Item(11, 15).Value = 0.103922044060685
$ws.Cells.Item(11, 16).Value = 0.1201207617372036
$ws.Cells.Item(11, 17).Value = 277.8265603771155
$ws.Cells.Item(11, 18).Value = 2500.43904339404
$ws.Cells.Item(11, 19).Value = 0.0497718439694778
$ws.Cells.Item(11, 20).Value = 0.05930848258337734

# Row 12
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 350.3919066666667
$ws.Cells.Item(12, 8).Value = 1051.17572
$ws.Cells.Item(12, 9).Value = 0.4789344206933965
$ws.Cells.Item(12, 10).Value = 0.4937404802104949
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 0.879594
$ws.Cells.Item(12, 14).Value = 2.638782
$ws.Cells.Item(12, 15).Value = 0.1152843201245645
$ws.Cells.Item(12, 16).Value = 0.1332541182661091
$ws.Cells.Item(12, 17).Value = 308.20261875256
$ws.Cells.Item(12, 18).Value = 2773.82356877304
$ws.Cells.Item(12, 19).Value = 0.05521362907389039
$ws.Cells.Item(12, 20).Value = 0.06579295234273479

# Row 13
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 350.3919066666667
$ws.Cells.Item(13, 8).Value = 1051.17572
$ws.Cells.Item(13, 9).Value = 0.4789344206933965
$ws.Cells.Item(13, 10).Value = 0.4937404802104949
$ws.Cells.Item(13, 11).Value = 2
$ws.Cells.Item(13, 13).Value = 1.8364275
$ws.Cells.Item(13, 14).Value = 3.672855
$ws.Cells.Item(13, 15).Value = 0.2406920645156217
$ws.Cells.Item(13, 16).Value = 0.1854730912005122
$ws.Cells.Item(13, 17).Value = 643.4693331801
$ws.Cells.Item(13, 18).Value = 3860.8159990806
$ws.Cells.Item(13, 19).Value = 0.1152757144842869
$ws.Cells.Item(13, 20).Value = 0.09157557311546582

# Row 14
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 243.8287033333334
$ws.Cells.Item(14, 8).Value = 731.4861100000001
$ws.Cells.Item(14, 9).Value = 0.3332781281688242
$ws.Cells.Item(14, 10).Value = 0.3435812836494235
$ws.Cells.Item(14, 11).Value = 2
$ws.Cells.Item(14, 13).Value = 1.2502825
$ws.Cells.Item(14, 14).Value = 2.500565
$ws.Cells.Item(14, 15).Value = 0.1638687485091313
$ws.Cells.Item(14, 16).Value = 0.1262743888059313
$ws.Cells.Item(14, 17).Value = 304.8547607753583
$ws.Cells.Item(14, 18).Value = 1829.12856465215
$ws.Cells.Item(14, 19).Value = 0.0546138697684911
$ws.Cells.Item(14, 20).Value = 0.04338551659798828

# Row 15
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 243.8287033333334
$ws.Cells.Item(15, 8).Value = 731.4861100000001
$ws.Cells.Item(15, 9).Value = 0.3332781281688242
$ws.Cells.Item(15, 10).Value = 0.3435812836494235
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 2.596814
$ws.Cells.Item(15, 14).Value = 7.790442000000001
$ws.Cells.Item(15, 15).Value = 0.3403524085884521
$ws.Cells.Item(15, 16).Value = 0.3934044114342388
$ws.Cells.Item(15, 17).Value = 633.1777904178467
$ws.Cells.Item(15, 18).Value = 5698.600113760621
$ws.Cells.Item(15, 19).Value = 0.1134320136521102
$ws.Cells.Item(15, 20).Value = 0.1351663926739217

# Row 16
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 243.8287033333334
$ws.Cells.Item(16, 8).Value = 731.4861100000001
$ws.Cells.Item(16, 9).Value = 0.3332781281688242
$ws.Cells.Item(16, 10).Value = 0.3435812836494235
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 0.2737596666666667
$ws.Cells.Item(16, 14).Value = 0.821279
$ws.Cells.Item(16, 15).Value = 0.03588041420154535
$ws.Cells.Item(16, 16).Value = 0.04147322855600493
$ws.Cells.Item(16, 17).Value = 66.75046454829889
$ws.Cells.Item(16, 18).Value = 600.75418093469
$ws.Cells.Item(16, 19).Value = 0.01195815728301313
$ws.Cells.Item(16, 20).Value = 0.0142494251043581

# Row 17
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 243.8287033333334
$ws.Cells.Item(17, 8).Value = 731.4861100000001
$ws.Cells.Item(17, 9).Value = 0.3332781281688242
$ws.Cells.Item(17, 10).Value = 0.3435812836494235
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 0.7929023333333333
$ws.Cells.Item(17, 14).Value = 2.378707
$ws.Cells.Item(17, 15).Value = 0.103922044060685
$ws.Cells.Item(17, 16).Value = 0.1201207617372036
$ws.Cells.Item(17, 17).Value = 193.3323478066411
$ws.Cells.Item(17, 18).Value = 1739.99113025977
$ws.Cells.Item(17, 19).Value = 0.03463494432002318
$ws.Cells.Item(17, 20).Value = 0.04127124551061496

# Row 18
$ws.Cells.Item(18, 5).Value = 3
$ws.Cells.Item(18, 7).Value = 243.8287033333334
$ws.Cells.Item(18, 8).Value = 731.4861100000001
$ws.Cells.Item(18, 9).Value = 0.3332781281688242
$ws.Cells.Item(18, 10).Value = 0.3435812836494235
$ws.Cells.Item(18, 11).Value = 3
$ws.Cells.Item(18, 13).Value = 0.879594
$ws.Cells.Item(18, 14).Value = 2.638782
$ws.Cells.Item(18, 15).Value = 0.1152843201245645
$ws.Cells.Item(18, 16).Value = 0.1332541182661091
$ws.Cells.Item(18, 17).Value = 214.47026447978
$ws.Cells.Item(18, 18).Value = 1930.23238031802
$ws.Cells.Item(18, 19).Value = 0.03842174241833039
$ws.Cells.Item(18, 20).Value = 0.04578362100544184

# Row 19
$ws.Cells.Item(19, 5).Value = 3
$ws.Cells.Item(19, 7).Value = 243.8287033333334
$ws.Cells.Item(19, 8).Value = 731.4861100000001
$ws.Cells.Item(19, 9).Value = 0.3332781281688242
$ws.Cells.Item(19, 10).Value = 0.3435812836494235
$ws.Cells.Item(19, 11).Value = 2
$ws.Cells.Item(19, 13).Value = 1.8364275
$ws.Cells.Item(19, 14).Value = 3.672855
$ws.Cells.Item(19, 15).Value = 0.2406920645156217
$ws.Cells.Item(19, 16).Value = 0.1854730912005122
$ws.Cells.Item(19, 17).Value = 447.773736090675
$ws.Cells.Item(19, 18).Value = 2686.64241654405
$ws.Cells.Item(19, 19).Value = 0.08021740072685626
$ws.Cells.Item(19, 20).Value = 0.06372508275709858

# Row 20
$ws.Cells.Item(20, 5).Value = 3
$ws.Cells.Item(20, 7).Value = 48.737294
$ws.Cells.Item(20, 8).Value = 146.211882
$ws.Cells.Item(20, 9).Value = 0.06661674320651284
$ws.Cells.Item(20, 10).Value = 0.06867617226847689
$ws.Cells.Item(20, 11).Value = 2
$ws.Cells.Item(20, 13).Value = 1.2502825
$ws.Cells.Item(20, 14).Value = 2.500565
$ws.Cells.Item(20, 15).Value = 0.1638687485091313
$ws.Cells.Item(20, 16).Value = 0.1262743888059313
$ws.Cells.Item(20, 17).Value = 60.935385785555
$ws.Cells.Item(20, 18).Value = 365.61231471333
$ws.Cells.Item(20, 19).Value = 0.01091640233900544
$ws.Cells.Item(20, 20).Value = 0.00867204167873277

# Row 21
$ws.Cells.Item(21, 5).Value = 3
$ws.Cells.Item(21, 7).Value = 48.737294
$ws.Cells.Item(21, 8).Value = 146.211882
$ws.Cells.Item(21, 9).Value = 0.06661674320651284
$ws.Cells.Item(21, 10).Value = 0.06867617226847689
$ws.Cells.Item(21, 11).Value = 3
$ws.Cells.Item(21, 13).Value = 2.596814
$ws.Cells.Item(21, 14).Value = 7.790442000000001
$ws.Cells.Item(21, 15).Value = 0.3403524085884521
$ws.Cells.Item(21, 16).Value = 0.3934044114342388
$ws.Cells.Item(21, 17).Value = 126.561687381316
$ws.Cells.Item(21, 18).Value = 1139.055186431844
$ws.Cells.Item(21, 19).Value = 0.02267316900265505
$ws.Cells.Item(21, 20).Value = 0.02701750913083655

# Row 22
$ws.Cells.Item(22, 5).Value = 3
$ws.Cells.Item(22, 7).Value = 48.737294
$ws.Cells.Item(22, 8).Value = 146.211882
$ws.Cells.Item(22, 9).Value = 0.06661674320651284
$ws.Cells.Item(22, 10).Value = 0.06867617226847689
$ws.Cells.Item(22, 11).Value = 3
$ws.Cells.Item(22, 13).Value = 0.2737596666666667
$ws.Cells.Item(22, 14).Value = 0.821279
$ws.Cells.Item(22, 15).Value = 0.03588041420154535
$ws.Cells.Item(22, 16).Value = 0.04147322855600493
$ws.Cells.Item(22, 17).Value = 13.34230535967533
$ws.Cells.Item(22, 18).Value = 120.080748237078
$ws.Cells.Item(22, 19).Value = 0.002390236339007663
$ws.Cells.Item(22, 20).Value = 0.00284822258884211

# Row 23
$ws.Cells.Item(23, 5).Value = 3
$ws.Cells.Item(23, 7).Value = 48.737294
$ws.Cells.Item(23, 8).Value = 146.211882
$ws.Cells.Item(23, 9).Value = 0.06661674320651284
$ws.Cells.Item(23, 10).Value = 0.06867617226847689
$ws.Cells.Item(23, 11).Value = 3
$ws.Cells.Item(23, 13).Value = 0.7929023333333333
$ws.Cells.Item(23, 14).Value = 2.378707
$ws.Cells.Item(23, 15).Value = 0.103922044060685
$ws.Cells.Item(23, 16).Value = 0.1201207617372036
$ws.Cells.Item(23, 17).Value = 38.64391413295266
$ws.Cells.Item(23, 18).Value = 347.795227196574
$ws.Cells.Item(23, 19).Value = 0.006922948122686566
$ws.Cells.Item(23, 20).Value = 0.00824943412608486

# Row 24
$ws.Cells.Item(24, 5).Value = 3
$ws.Cells.Item(24, 7).Value = 48.737294
$ws.Cells.Item(24, 8).Value = 146.211882
$ws.Cells.Item(24, 9).Value = 0.06661674320651284
$ws.Cells.Item(24, 10).Value = 0.06867617226847689
$ws.Cells.Item(24, 11).Value = 3
$ws.Cells.Item(24, 13).Value = 0.879594
$ws.Cells.Item(24, 14).Value = 2.638782
$ws.Cells.Item(24, 15).Value = 0.1152843201245645
$ws.Cells.Item(24, 16).Value = 0.1332541182661091
$ws.Cells.Item(24, 17).Value = 42.869031378636
$ws.Cells.Item(24, 18).Value = 385.821282407724
$ws.Cells.Item(24, 19).Value = 0.007679865949475536
$ws.Cells.Item(24, 20).Value = 0.009151382781527301

# Row 25
$ws.Cells.Item(25, 5).Value = 3
$ws.Cells.Item(25, 7).Value = 48.737294
$ws.Cells.Item(25, 8).Value = 146.211882
$ws.Cells.Item(25, 9).Value = 0.06661674320651284
$ws.Cells.Item(25, 10).Value = 0.06867617226847689
$ws.Cells.Item(25, 11).Value = 2
$ws.Cells.Item(25, 13).Value = 1.8364275
$ws.Cells.Item(25, 14).Value = 3.672855
$ws.Cells.Item(25, 15).Value = 0.2406920645156217
$ws.Cells.Item(25, 16).Value = 0.1854730912005122
$ws.Cells.Item(25, 17).Value = 89.50250697718499
$ws.Cells.Item(25, 18).Value = 537.01504186311
$ws.Cells.Item(25, 19).Value = 0.01603412145368259
$ws.Cells.Item(25, 20).Value = 0.0127375819624533

# Row 26
$ws.Cells.Item(26, 5).Value = 3
$ws.Cells.Item(26, 7).Value = 22.832077
$ws.Cells.Item(26, 8).Value = 68.49623099999999
$ws.Cells.Item(26, 9).Value = 0.03120810544755168
$ws.Cells.Item(26, 10).Value = 0.03217289111905
$ws.Cells.Item(26, 11).Value = 2
$ws.Cells.Item(26, 13).Value = 1.2502825
$ws.Cells.Item(26, 14).Value = 2.500565
$ws.Cells.Item(26, 15).Value = 0.1638687485091313
$ws.Cells.Item(26, 16).Value = 0.1262743888059313
$ws.Cells.Item(26, 17).Value = 28.5465463117525
$ws.Cells.Item(26, 18).Value = 171.279277870515
$ws.Cells.Item(26, 19).Value = 0.005114033183031298
$ws.Cells.Item(26, 20).Value = 0.004062612162177815

# Row 27
$ws.Cells.Item(27, 5).Value = 3
$ws.Cells.Item(27, 7).Value = 22.832077
$ws.Cells.Item(27, 8).Value = 68.49623099999999
$ws.Cells.Item(27, 9).Value = 0.03120810544755168
$ws.Cells.Item(27, 10).Value = 0.03217289111905
$ws.Cells.Item(27, 11).Value = 3
$ws.Cells.Item(27, 13).Value = 2.596814
$ws.Cells.Item(27, 14).Value = 7.790442000000001
$ws.Cells.Item(27, 15).Value = 0.3403524085884521
$ws.Cells.Item(27, 16).Value = 0.3934044114342388
$ws.Cells.Item(27, 17).Value = 59.290657202678
$ws.Cells.Item(27, 18).Value = 533.615914824102
$ws.Cells.Item(27, 19).Value = 0.01062175385655661
$ws.Cells.Item(27, 20).Value = 0.01265695729482771

# Row 28
$ws.Cells.Item(28, 5).Value = 3
$ws.Cells.Item(28, 7).Value = 22.832077
$ws.Cells.Item(28, 8).Value = 68.49623099999999
$ws.Cells.Item(28, 9).Value = 0.03120810544755168
$ws.Cells.Item(28, 10).Value = 0.03217289111905
$ws.Cells.Item(28, 11).Value = 3
$ws.Cells.Item(28, 13).Value = 0.2737596666666667
$ws.Cells.Item(28, 14).Value = 0.821279
$ws.Cells.Item(28, 15).Value = 0.03588041420154535
$ws.Cells.Item(28, 16).Value = 0.04147322855600493
$ws.Cells.Item(28, 17).Value = 6.250501788827666
$ws.Cells.Item(28, 18).Value = 56.25451609944899
$ws.Cells.Item(28, 19).Value = 0.001119759749903658
$ws.Cells.Item(28, 20).Value = 0.001334313666687822

# Row 29
$ws.Cells.Item(29, 5).Value = 3
$ws.Cells.Item(29, 7).Value = 22.832077
$ws.Cells.Item(29, 8).Value = 68.49623099999999
$ws.Cells.Item(29, 9).Value = 0.03120810544755168
$ws.Cells.Item(29, 10).Value = 0.03217289111905
$ws.Cells.Item(29, 11).Value = 3
$ws.Cells.Item(29, 13).Value = 0.7929023333333333
$ws.Cells.Item(29, 14).Value = 2.378707
$ws.Cells.Item(29, 15).Value = 0.103922044060685
$ws.Cells.Item(29, 16).Value = 0.1201207617372036
$ws.Cells.Item(29, 17).Value = 18.10360712814633
$ws.Cells.Item(29, 18).Value = 162.932464153317
$ws.Cells.Item(29, 19).Value = 0.00324321010937097
$ws.Cells.Item(29, 20).Value = 0.003864632188508398

# Row 30
$ws.Cells.Item(30, 5).Value = 3
$ws.Cells.Item(30, 7).Value = 22.832077
$ws.Cells.Item(30, 8).Value = 68.49623099999999
$ws.Cells.Item(30, 9).Value = 0.03120810544755168
$ws.Cells.Item(30, 10).Value = 0.03217289111905
$ws.Cells.Item(30, 11).Value = 3
$ws.Cells.Item(30, 13).Value = 0.879594
$ws.Cells.Item(30, 14).Value = 2.638782
$ws.Cells.Item(30, 15).Value = 0.1152843201245645
$ws.Cells.Item(30, 16).Value = 0.1332541182661091
$ws.Cells.Item(30, 17).Value = 20.082957936738
$ws.Cells.Item(30, 18).Value = 180.746621430642
$ws.Cells.Item(30, 19).Value = 0.003597805218896714
$ws.Cells.Item(30, 20).Value = 0.004287170238140539

# Row 31
$ws.Cells.Item(31, 5).Value = 3
$ws.Cells.Item(31, 7).Value = 22.832077
$ws.Cells.Item(31, 8).Value = 68.49623099999999
$ws.Cells.Item(31, 9).Value = 0.03120810544755168
$ws.Cells.Item(31, 10).Value = 0.03217289111905
$ws.Cells.Item(31, 11).Value = 2
$ws.Cells.Item(31, 13).Value = 1.8364275
$ws.Cells.Item(31, 14).Value = 3.672855
$ws.Cells.Item(31, 15).Value = 0.2406920645156217
$ws.Cells.Item(31, 16).Value = 0.1854730912005122
$ws.Cells.Item(31, 17).Value = 41.9294540849175
$ws.Cells.Item(31, 18).Value = 251.576724509505
$ws.Cells.Item(31, 19).Value = 0.007511543329792434
$ws.Cells.Item(31, 20).Value = 0.005967205568707711

# Row 32
$ws.Cells.Item(32, 5).Value = 2
$ws.Cells.Item(32, 7).Value = 25.082339
$ws.Cells.Item(32, 8).Value = 50.164678
$ws.Cells.Item(32, 9).Value = 0.03428388404538221
$ws.Cells.Item(32, 10).Value = 0.02356250409334498
$ws.Cells.Item(32, 11).Value = 2
$ws.Cells.Item(32, 13).Value = 1.2502825
$ws.Cells.Item(32, 14).Value = 2.500565
$ws.Cells.Item(32, 15).Value = 0.1638687485091313
$ws.Cells.Item(32, 16).Value = 0.1262743888059313
$ws.Cells.Item(32, 17).Value = 31.3600095107675
$ws.Cells.Item(32, 18).Value = 125.44003804307
$ws.Cells.Item(32, 19).Value = 0.005618057172548958
$ws.Cells.Item(32, 20).Value = 0.002975340803124392

# Row 33
$ws.Cells.Item(33, 5).Value = 2
$ws.Cells.Item(33, 7).Value = 25.082339
$ws.Cells.Item(33, 8).Value = 50.164678
$ws.Cells.Item(33, 9).Value = 0.03428388404538221
$ws.Cells.Item(33, 10).Value = 0.02356250409334498
$ws.Cells.Item(33, 11).Value = 3
$ws.Cells.Item(33, 13).Value = 2.596814
$ws.Cells.Item(33, 14).Value = 7.790442000000001
$ws.Cells.Item(33, 15).Value = 0.3403524085884521
$ws.Cells.Item(33, 16).Value = 0.3934044114342388
$ws.Cells.Item(33, 17).Value = 65.134169067946
$ws.Cells.Item(33, 18).Value = 390.805014407676
$ws.Cells.Item(33, 19).Value = 0.01166860251061304
$ws.Cells.Item(33, 20).Value = 0.009269593054759223

# Row 34
$ws.Cells.Item(34, 5).Value = 2
$ws.Cells.Item(34, 7).Value = 25.082339
$ws.Cells.Item(34, 8).Value = 50.164678
$ws.Cells.Item(34, 9).Value = 0.03428388404538221
$ws.Cells.Item(34, 10).Value = 0.02356250409334498
$ws.Cells.Item(34, 11).Value = 3
$ws.Cells.Item(34, 13).Value = 0.2737596666666667
$ws.Cells.Item(34, 14).Value = 0.821279
$ws.Cells.Item(34, 15).Value = 0.03588041420154535
$ws.Cells.Item(34, 16).Value = 0.04147322855600493
$ws.Cells.Item(34, 17).Value = 6.866532763860334
$ws.Cells.Item(34, 18).Value = 41.199196583162
$ws.Cells.Item(34, 19).Value = 0.001230119959986066
$ws.Cells.Item(34, 20).Value = 0.000977213117615098

# Row 35
$ws.Cells.Item(35, 5).Value = 2
$ws.Cells.Item(35, 7).Value = 25.082339
$ws.Cells.Item(35, 8).Value = 50.164678
$ws.Cells.Item(35, 9).Value = 0.03428388404538221
$ws.Cells.Item(35, 10).Value = 0.02356250409334498
$ws.Cells.Item(35, 11).Value = 3
$ws.Cells.Item(35, 13).Value = 0.7929023333333333
$ws.Cells.Item(35, 14).Value = 2.378707
$ws.Cells.Item(35, 15).Value = 0.103922044060685
$ws.Cells.Item(35, 16).Value = 0.1201207617372036
$ws.Cells.Item(35, 17).Value = 19.88784511855767
$ws.Cells.Item(35, 18).Value = 119.327070711346
$ws.Cells.Item(35, 19).Value = 0.003562851308335626
$ws.Cells.Item(35, 20).Value = 0.002830345940128576

# Row 36
$ws.Cells.Item(36, 5).Value = 2
$ws.Cells.Item(36, 7).Value = 25.082339
$ws.Cells.Item(36, 8).Value = 50.164678
$ws.Cells.Item(36, 9).Value = 0.03428388404538221
$ws.Cells.Item(36, 10).Value = 0.02356250409334498
$ws.Cells.Item(36, 11).Value = 3
$ws.Cells.Item(36, 13).Value = 0.879594
$ws.Cells.Item(36, 14).Value = 2.638782
$ws.Cells.Item(36, 15).Value = 0.1152843201245645
$ws.Cells.Item(36, 16).Value = 0.1332541182661091
$ws.Cells.Item(36, 17).Value = 22.062274890366
$ws.Cells.Item(36, 18).Value = 132.373649342196
$ws.Cells.Item(36, 19).Value = 0.003952394263401294
$ws.Cells.Item(36, 20).Value = 0.003139800707100271

# Row 37
$ws.Cells.Item(37, 5).Value = 2
$ws.Cells.Item(37, 7).Value = 25.082339
$ws.Cells.Item(37, 8).Value = 50.164678
$ws.Cells.Item(37, 9).Value = 0.03428388404538221
$ws.Cells.Item(37, 10).Value = 0.02356250409334498
$ws.Cells.Item(37, 11).Value = 2
$ws.Cells.Item(37, 13).Value = 1.8364275
$ws.Cells.Item(37, 14).Value = 3.672855
$ws.Cells.Item(37, 15).Value = 0.2406920645156217
$ws.Cells.Item(37, 16).Value = 0.1854730912005122
$ws.Cells.Item(37, 17).Value = 46.0618971039225
$ws.Cells.Item(37, 18).Value = 184.24758841569
$ws.Cells.Item(37, 19).Value = 0.008251858830497227
$ws.Cells.Item(37, 20).Value = 0.004370210470617415
